$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ignore" modifier column (E) with header + first data row centered
$ws.Range("E1").Value = "Ignore"
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E2").Value = "Ignore"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108

# New "Feel" entries (column C)
$ws.Range("C14").Value = "Piano-house"
$ws.Range("C15").Value = "Roller"
$ws.Range("C16").Value = "Rock"
$ws.Range("C17").Value = "Metal"

# New "Modifier" entries (column B)
$ws.Range("B15").Value = "Jump-Up"
$ws.Range("B16").Value = "Melodic"

# New "DJ Tool" entry (column D)
$ws.Range("D7").Value = "Tool"

# New "Genre" entry (column A)
$ws.Range("A15").Value = "Ambient"

$ws.Range("A15").Select() | Out-Null
